# Auto-generated update script for Goblin_Profits market-price refresh
# Applies per-cell value updates (H:currentAveragePrice .. N:LeveProfitHQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-crafting tables.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = @(
    ,(64, 8, 7280.121)
    ,(64, 9, 3473.0833)
    ,(64, 10, 9455.571)
    ,(64, 11, 3473.0833)
    ,(64, 12, 9455.571)
    ,(64, 13, -3225.0833)
    ,(64, 14, -9951.571)
    ,(67, 8, 7280.121)
    ,(67, 9, 3473.0833)
    ,(67, 10, 9455.571)
    ,(67, 11, 3473.0833)
    ,(67, 12, 9455.571)
    ,(67, 13, -2615.0833)
    ,(67, 14, -11171.571)
    ,(80, 8, 1729.1578)
    ,(80, 10, 2619.4)
    ,(80, 12, 7858.200000000001)
    ,(80, 14, -9854.200000000001)
    ,(83, 8, 1729.1578)
    ,(83, 10, 2619.4)
    ,(83, 12, 23574.6)
    ,(83, 14, -33558.60000000001)
)
foreach ($u in $ALC_updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = @(
    ,(32, 8, 83539.125)
    ,(32, 9, 98774.45)
    ,(32, 10, 16249.75)
    ,(32, 11, 98774.45)
    ,(32, 12, 16249.75)
    ,(32, 13, -98487.45)
    ,(32, 14, -16823.75)
    ,(61, 8, 6039.162)
    ,(61, 9, 4808.143)
    ,(61, 11, 4808.143)
    ,(61, 13, -4596.143)
    ,(136, 8, 6039.162)
    ,(136, 9, 4808.143)
    ,(136, 11, 14424.429)
    ,(136, 13, -11874.429)
    ,(139, 8, 241666.67)
    ,(139, 10, 241666.67)
    ,(139, 12, 241666.67)
    ,(139, 14, -251946.67)
)
foreach ($u in $ARM_updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = @(
    ,(20, 8, 4486.9062)
    ,(20, 9, 5908.0454)
    ,(20, 10, 1360.4)
    ,(20, 11, 5908.0454)
    ,(20, 12, 1360.4)
    ,(20, 13, -5661.0454)
    ,(20, 14, -1854.4)
    ,(87, 8, 81500)
    ,(87, 10, 81500)
    ,(87, 12, 81500)
    ,(87, 14, -83996)
    ,(90, 8, 81500)
    ,(90, 10, 81500)
    ,(90, 12, 244500)
    ,(90, 14, -256980)
    ,(99, 8, 4318.1875)
    ,(99, 9, 4111.5)
    ,(99, 10, 4524.875)
    ,(99, 11, 4111.5)
    ,(99, 12, 4524.875)
    ,(99, 13, -2613.5)
    ,(99, 14, -7520.875)
    ,(105, 8, 12738.23)
    ,(105, 9, 17273.875)
    ,(105, 11, 17273.875)
    ,(105, 13, -15526.875)
    ,(107, 8, 3491.5208)
    ,(107, 9, 2555.5)
    ,(107, 10, 6299.5835)
    ,(107, 11, 2555.5)
    ,(107, 12, 6299.5835)
    ,(107, 13, -635.5)
    ,(107, 14, -10139.5835)
    ,(133, 8, 67316.664)
    ,(133, 10, 67316.664)
    ,(133, 12, 67316.664)
    ,(133, 14, -77436.664)
)
foreach ($u in $BSM_updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = @(
    ,(31, 8, 3851.577)
    ,(31, 9, 1967.9412)
    ,(31, 10, 7409.5557)
    ,(31, 11, 1967.9412)
    ,(31, 12, 7409.5557)
    ,(31, 13, -1672.9412)
    ,(31, 14, -7999.5557)
    ,(34, 8, 3851.577)
    ,(34, 9, 1967.9412)
    ,(34, 10, 7409.5557)
    ,(34, 11, 1967.9412)
    ,(34, 12, 7409.5557)
    ,(34, 13, -1765.9412)
    ,(34, 14, -7813.5557)
    ,(97, 8, 0)
    ,(97, 10, 0)
    ,(97, 12, 0)
    ,(97, 14, $null)
    ,(116, 8, 0)
    ,(116, 10, 0)
    ,(116, 12, 0)
    ,(116, 14, $null)
    ,(134, 8, 40550.41)
    ,(134, 9, 42522.312)
    ,(134, 10, 9000)
    ,(134, 11, 127566.936)
    ,(134, 12, 27000)
    ,(134, 13, -125031.936)
    ,(134, 14, -32070)
    ,(141, 8, 201558.94)
    ,(141, 10, 228721.84)
    ,(141, 12, 228721.84)
    ,(141, 14, -239081.84)
)
foreach ($u in $CRP_updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$CUL_updates = @(
    ,(51, 8, 840.2)
    ,(51, 9, 565.6667)
    ,(51, 11, 1697.0001)
    ,(51, 13, -1237.0001)
    ,(57, 8, 6999.778)
    ,(57, 9, 6999)
    ,(57, 10, 6999.875)
    ,(57, 11, 20997)
    ,(57, 12, 20999.625)
    ,(57, 13, -20438)
    ,(57, 14, -22117.625)
    ,(63, 8, 2629.3333)
    ,(63, 9, 944)
    ,(63, 11, 2832)
    ,(63, 13, -2083)
    ,(66, 8, 2629.3333)
    ,(66, 9, 944)
    ,(66, 11, 8496)
    ,(66, 13, -4752)
    ,(68, 8, 3408.5)
    ,(68, 9, 3738)
    ,(68, 10, 3225.4443)
    ,(68, 11, 11214)
    ,(68, 12, 9676.332900000001)
    ,(68, 13, -10403)
    ,(68, 14, -11298.3329)
    ,(71, 8, 3408.5)
    ,(71, 9, 3738)
    ,(71, 10, 3225.4443)
    ,(71, 11, 33642)
    ,(71, 12, 29028.9987)
    ,(71, 13, -29586)
    ,(71, 14, -37140.9987)
    ,(82, 8, 59328.89)
    ,(82, 9, 59500)
    ,(82, 11, 178500)
    ,(82, 13, -178094)
    ,(85, 8, 59328.89)
    ,(85, 9, 59500)
    ,(85, 11, 178500)
    ,(85, 13, -177096)
    ,(114, 8, 1228)
    ,(114, 9, 1425)
    ,(114, 10, 1031)
    ,(114, 11, 4275)
    ,(114, 12, 3093)
    ,(114, 13, -1021)
    ,(114, 14, -9601)
    ,(115, 8, 7495.6924)
    ,(115, 9, 3636.125)
    ,(115, 10, 13671)
    ,(115, 11, 10908.375)
    ,(115, 12, 41013)
    ,(115, 13, -9733.375)
    ,(115, 14, -43363)
    ,(119, 8, 668.4286)
    ,(119, 9, 596.5)
    ,(119, 11, 1789.5)
    ,(119, 13, 3048.5)
    ,(123, 8, 1696.4)
    ,(123, 9, 1120.5)
    ,(123, 11, 3361.5)
    ,(123, 13, -911.5)
    ,(132, 8, 2299.2415)
    ,(132, 9, 2346)
    ,(132, 10, 2295.7778)
    ,(132, 11, 21114)
    ,(132, 12, 20662.0002)
    ,(132, 13, -18584)
    ,(132, 14, -25722.0002)
    ,(138, 8, 2399.9167)
    ,(138, 9, 2260.889)
    ,(138, 10, 2483.3333)
    ,(138, 11, 6782.667)
    ,(138, 12, 7449.999899999999)
    ,(138, 13, -1642.667)
    ,(138, 14, -17729.9999)
)
foreach ($u in $CUL_updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = @(
    ,(92, 8, 15583.167)
    ,(92, 10, 15583.167)
    ,(92, 12, 15583.167)
    ,(92, 14, -19327.167)
    ,(94, 8, 0)
    ,(94, 10, 0)
    ,(94, 12, 0)
    ,(94, 14, $null)
    ,(102, 8, 14287173)
    ,(102, 9, 17858000)
    ,(102, 10, 3869)
    ,(102, 11, 17858000)
    ,(102, 12, 3869)
    ,(102, 13, -17856378)
    ,(102, 14, -7113)
)
foreach ($u in $GSM_updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = @(
    ,(46, 8, 1271.2188)
    ,(46, 9, 1006.2941)
    ,(46, 10, 1571.4667)
    ,(46, 11, 1006.2941)
    ,(46, 12, 1571.4667)
    ,(46, 13, -818.2941)
    ,(46, 14, -1947.4667)
    ,(106, 8, 27562)
    ,(106, 10, 24356.572)
    ,(106, 12, 24356.572)
    ,(106, 14, -26880.572)
    ,(122, 8, 4988.1113)
    ,(122, 9, 5216.96)
    ,(122, 10, 4468)
    ,(122, 11, 15650.88)
    ,(122, 12, 13404)
    ,(122, 13, -13200.88)
    ,(122, 14, -18304)
)
foreach ($u in $LTW_updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = @(
    ,(62, 8, 10589.728)
    ,(62, 10, 12916.583)
    ,(62, 12, 12916.583)
    ,(62, 14, -14164.583)
    ,(65, 8, 10589.728)
    ,(65, 10, 12916.583)
    ,(65, 12, 64582.915)
    ,(65, 14, -70822.91500000001)
    ,(107, 8, 650.5454999999999)
    ,(107, 9, 665.2857)
    ,(107, 11, 1995.8571)
    ,(107, 13, -75.85710000000017)
    ,(132, 8, 3259.9111)
    ,(132, 9, 2514.973)
    ,(132, 10, 6705.25)
    ,(132, 11, 7544.919)
    ,(132, 12, 20115.75)
    ,(132, 13, -5014.919)
    ,(132, 14, -25175.75)
)
foreach ($u in $WVR_updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
